# MCHP Football Fund - record June,18 (column F) payments of 300 each for
# Ganesh Kumar (row 8), Vignesh (row 28), Vikram (row 29) and Yogesh (row 33).
# Yogesh's previous (wrong) 100 in F33 becomes 300, and Vignesh's previous
# 100 that had been recorded in F28 (June) is corrected to live in G28
# (July,18) instead, with June (F28) getting its own 300 entry.
#
# Commit message: "yogesh, vikram, vignesh, ganesh paid"
#
# All of the "Total" row (35) and running "Balance" column (E40:E52) cells
# are formula-driven and recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ganesh Kumar (row 8) paid 300 for June,18
$ws.Range("F8").Value = 300

# Vignesh (row 28): 300 for June,18, and his existing 100 payment is
# (re)recorded under July,18 instead of June,18
$ws.Range("F28").Value = 300
$ws.Range("G28").Value = 100

# Vikram (row 29) paid 300 for June,18
$ws.Range("F29").Value = 300

# Yogesh (row 33): June,18 payment corrected from 100 to 300
$ws.Range("F33").Value = 300

# Reflect the author's resulting view state (scrolled down, F29 selected)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F29").Select()
